$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) from 45737 to 45738
# for every existing data row (2 through 44).
for ($r = 2; $r -le 44; $r++) {
    $ws.Cells.Item($r, 3).Value = 45738
}

# Row 44 picks up an explicit row height in the target file; match it.
$ws.Rows.Item(44).RowHeight = 15

# Append the new data row (45) for case A 13183-2025.
$ws.Cells.Item(45, 1).Value = "A 13183-2025"

$ws.Cells.Item(45, 2).Value = 45735
$ws.Cells.Item(45, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(45, 3).Value = 45738
$ws.Cells.Item(45, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(45, 4).Value = "OKÄNT"
$ws.Cells.Item(45, 5).Value = "OKÄNT"

$ws.Cells.Item(45, 7).Value = 0.6
$ws.Cells.Item(45, 8).Value = 0
$ws.Cells.Item(45, 9).Value = 0
$ws.Cells.Item(45, 10).Value = 0
$ws.Cells.Item(45, 11).Value = 0
$ws.Cells.Item(45, 12).Value = 0
$ws.Cells.Item(45, 13).Value = 0
$ws.Cells.Item(45, 14).Value = 0
$ws.Cells.Item(45, 15).Value = 0
$ws.Cells.Item(45, 16).Value = 0
$ws.Cells.Item(45, 17).Value = 0

# R45 is a blank, wrap-text-styled cell (matches the other rows' "Artnamn" cells).
$ws.Range("R45").WrapText = $true
